$d = $word.ActiveDocument

# Locate the "Gender: {{gender}}" paragraph; the new "CRN: {{crn}}" paragraph
# must be inserted directly after it (and before the "CRO No: {{cro_number}}"
# paragraph that currently follows it).
$genderPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{{gender}}*") {
        $genderPara = $p
        break
    }
}

if ($genderPara -ne $null) {
    $followingPara = $genderPara.Next()
    $followRange = $followingPara.Range
    $followRange.Collapse(1)
    # Inserting a paragraph break before the following paragraph's range
    # clones that paragraph's formatting (borders/tabs/spacing/rPr) onto
    # the freshly created empty paragraph - matching the target pPr.
    $followRange.InsertParagraphBefore()
}

# Re-locate the Gender paragraph (collection was mutated) and grab the new
# empty paragraph that now follows it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{{gender}}*") {
        $newPara = $p.Next()
        $pr = $newPara.Range
        $pr.Collapse(1)
        $startPos = $pr.Start

        # Insert the full run of text first (inherits surrounding formatting),
        # then re-style only the "CRN:" label as bold + purple, leaving the
        # placeholder text in the paragraph's plain (non-bold) style.
        $pr.InsertAfter("CRN: {{crn}}")

        $labelRange = $d.Range($startPos, $startPos + 4)
        $labelRange.Font.Bold = $true
        $labelRange.Font.Color = 8388736
        break
    }
}
